# PM21 Tidsregistrering for Lukas.xlsx - add a new logged time entry (row 7)
# on sheet "Ark1": a Test Analyst task "Lav SSD02 - ..." on 2020-03-09,
# 13:00-16:30 (3.5 hours), with an estimate of 60 minutes. Also refresh the
# sheet view (zoom back to 100%, move the selection to F10).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row of time-tracking data
$ws.Range("A7").Value = 'Lav SSD02 - Fundet og rettet flere fejl til vedhørende "02" filer'
$ws.Range("B7").Value = 'Test Analyst'
$ws.Range("C7").Value = 43899
$ws.Range("D7").Value = 0.54166666666666663
$ws.Range("E7").Value = 0.6875
$ws.Range("F7").Value = 60

# Recalculate so the running-total column (H) picks up the new row
$excel.Calculate()

# Restore the view: zoom to 100% and move the selection to F10
$excel.ActiveWindow.Zoom = 100
$ws.Range("F10").Select() | Out-Null
